$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.697.41"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.798.93"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.81"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.18"
$ws.Range("E6").Value = "  +1.03%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  +0.87%  "

$ws.Range("E9").Value = "  +1.80%  "

$ws.Range("E10").Value = "  -0.92%  "

$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("E12").Value = "  +0.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.88"
$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.440.25"
$ws.Range("E14").Value = "  +1.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.859.93"
$ws.Range("E15").Value = "  +2.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.56"
$ws.Range("E16").Value = "  +4.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.654.38"
$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("E18").Value = "  +2.54%  "

$ws.Range("E19").Value = "  +0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.56"
$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.90"
$ws.Range("E21").Value = "  -5.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.702"
$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("E23").Value = "  +2.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.44"
$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.10"
$ws.Range("E25").Value = "  +2.68%  "

$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.01"
$ws.Range("E28").Value = "  +0.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.939.99"
$ws.Range("E29").Value = "  +1.14%  "

$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("E31").Value = "  +2.93%  "

$ws.Range("E32").Value = "  +1.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.58"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("E35").Value = "  -0.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0999"
$ws.Range("E36").Value = "  +0.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.38"
$ws.Range("E37").Value = "  +2.94%  "

$ws.Range("E38").Value = "  +0.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.997"
$ws.Range("E39").Value = "  +1.15%  "

$ws.Range("E40").Value = "  +0.92%  "

$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.14"
$ws.Range("E43").Value = "  +3.16%  "

$ws.Range("E44").Value = "  +1.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.22"
$ws.Range("E45").Value = "  -1.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.33"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.00"
$ws.Range("E47").Value = "  +0.62%  "

$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.36"
$ws.Range("E48").Value = "  +12.24%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "394.65"
$ws.Range("E49").Value = "  +2.15%  "

$ws.Range("E50").Value = "  +1.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.85"
$ws.Range("E51").Value = "  +7.37%  "
